# Refresh crypto price/volume snapshot (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "51.949.49"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "'" + "2.790.51"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'" + "361.76"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'" + "109.81"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "'" + "40.09"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").Value = "'" + "0.0847"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'" + "19.49"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").Value = "'" + "7.57"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "'" + "3.226.55"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "'" + "2.779.03"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").Value = "'" + "0.939"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("D18").Value = "'" + "51.897.17"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'" + "7.51"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "'" + "13.12"
$ws.Range("E21").Value = "  -3.49%  "
$ws.Range("D22").Value = "'" + "0.0₃0976"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "'" + "70.35"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'" + "269.63"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "'" + "2.76"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").Value = "'" + "26.54"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'" + "0.160"
$ws.Range("E28").Value = "  +14.19%  "
$ws.Range("D29").Value = "'" + "10.30"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "'" + "2.21"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "'" + "0.0474"
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "'" + "51.89"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").Value = "'" + "34.03"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'" + "5.74"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "'" + "0.0844"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'" + "5.24"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'" + "18.94"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").Value = "'" + "3.22"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").Value = "'" + "2.58"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "'" + "2.25"
$ws.Range("D44").Value = "'" + "119.81"
$ws.Range("E44").Value = "  -6.47%  "
$ws.Range("D45").Value = "'" + "21.95"
$ws.Range("E45").Value = "  -8.48%  "
$ws.Range("D46").Value = "'" + "2.087.94"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "'" + "3.26"
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'" + "0.958"
$ws.Range("E50").Value = "  -4.83%  "
$ws.Range("D51").Value = "'" + "8.88"
$ws.Range("E51").Value = "  -1.80%  "
